$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c5 = [char]0x2085
$c8 = [char]0x2088

$ws.Range("D2").Value = "26.191.02"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "1.659.20"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'217.80"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'0.5209"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'0.2661"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").Value = "'0.06299"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'20.98"
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "1.659.52"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "'4.418"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "1.884.80"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "'0.5447"
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").Value = "0.0${c5}8181"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").Value = "'64.61"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "26.213.92"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'4.651"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "'192.05"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").Value = "'10.12"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("D23").Value = "'6.066"
$ws.Range("E23").Value = "  -4.87%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'138.49"
$ws.Range("E25").Value = "  -4.06%  "
$ws.Range("D26").Value = "'0.1235"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").Value = "'7.195"
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("D28").Value = "'16.13"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'1.412"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'0.05977"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("D31").Value = "'1.281"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'3.575"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'3.309"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").Value = "'1.634"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").Value = "'0.9758"
$ws.Range("E35").Value = "  -3.68%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.782"
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.410"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "'0.5852"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("D39").Value = "'0.01583"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").Value = "'5.931"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").Value = "'0.8611"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "1.031.89"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").Value = "'99.48"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "1.800.23"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'56.90"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0${c8}107"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").Value = "'8.062"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "'0.05183"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").Value = "'0.4230"
$ws.Range("E51").Value = "  -0.26%  "
